# Auto-generated edit script applying the cryptos.xlsx GitHub Actions price/volume update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.875.01"
$ws.Range("E2").Value = "  -0.03%  "
Set-TextValue "D3" "1.887.56"
$ws.Range("E3").Value = "  -0.39%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "0.7691"
$ws.Range("E5").Value = "  -0.89%  "
Set-TextValue "D6" "242.64"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.04%  "
Set-TextValue "D9" "25.63"
$ws.Range("E9").Value = "  +0.37%  "
Set-TextValue "D10" "0.07165"
$ws.Range("E10").Value = "  -5.26%  "
Set-TextValue "D11" "0.08577"
$ws.Range("E11").Value = "  +5.71%  "
Set-TextValue "D12" "0.7634"
$ws.Range("E12").Value = "  -1.04%  "
Set-TextValue "D13" "1.917.48"
$ws.Range("E13").Value = "  +2.20%  "
Set-TextValue "D14" "5.363"
$ws.Range("E14").Value = "  -2.21%  "
Set-TextValue "D15" "93.57"
$ws.Range("E15").Value = "  +1.30%  "
Set-TextValue "D16" "6.145"
$ws.Range("E16").Value = "  -1.38%  "
Set-TextValue "D17" "29.911.97"
$ws.Range("E17").Value = "  +0.41%  "
Set-TextValue "D18" "13.76"
$ws.Range("E18").Value = "  -1.87%  "
Set-TextValue "D19" "244.47"
$ws.Range("E19").Value = "  -0.08%  "
Set-TextValue "D20" "0.000007801"
$ws.Range("E20").Value = "  -1.50%  "
Set-TextValue "D21" "2.170.57"
$ws.Range("E21").Value = "  +3.77%  "
Set-TextValue "D22" "0.9995"
$ws.Range("E22").Value = "  -0.03%  "
Set-TextValue "D23" "8.010"
$ws.Range("E23").Value = "  -1.14%  "
Set-TextValue "D24" "1.001"
Set-TextValue "D25" "0.1636"
$ws.Range("E25").Value = "  +4.18%  "
Set-TextValue "D26" "9.378"
$ws.Range("E26").Value = "  -0.92%  "
Set-TextValue "D27" "162.72"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -0.51%  "
Set-TextValue "D29" "2.032"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("E30").Value = "  +1.81%  "
Set-TextValue "D31" "1.540"
$ws.Range("E31").Value = "  -0.73%  "
Set-TextValue "D32" "4.516"
$ws.Range("E32").Value = "  +0.63%  "
Set-TextValue "D33" "4.094"
Set-TextValue "D34" "0.05458"
$ws.Range("E34").Value = "  -0.98%  "
Set-TextValue "D35" "1.239"
$ws.Range("E35").Value = "  -1.85%  "
Set-TextValue "D36" "0.7427"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("E37").Value = "  +0.45%  "
Set-TextValue "D38" "2.696"
$ws.Range("E38").Value = "  +2.05%  "
Set-TextValue "D39" "0.01953"
$ws.Range("E39").Value = "  +1.30%  "
Set-TextValue "D40" "2.783"
$ws.Range("E40").Value = "  -0.29%  "
Set-TextValue "D41" "0.4464"
$ws.Range("E41").Value = "  +0.25%  "
Set-TextValue "D42" "1.109.43"
$ws.Range("E42").Value = "  -4.58%  "
Set-TextValue "D43" "6.079"
$ws.Range("E43").Value = "  +2.31%  "
Set-TextValue "D44" "73.07"
$ws.Range("E44").Value = "  -1.43%  "
Set-TextValue "D45" "0.8520"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  -0.01%  "
Set-TextValue "D47" "102.41"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D48" "7.636"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.860"
$ws.Range("E49").Value = "  -2.14%  "
Set-TextValue "D50" "3.007"
$ws.Range("E50").Value = "  -3.96%  "
Set-TextValue "D51" "2.054.34"
$ws.Range("E51").Value = "  +1.33%  "
